# Weekly update: insert a new price record as row 17 (new week's reading),
# pushing the previously-existing rows 17-38 down to 18-39 unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 17, shifting rows 17:38 down
# to 18:39 (formatting of the row above, incl. the date style on column D,
# carries down automatically - matching Excel's native Insert behaviour).
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the new weekly record.
$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value = "La Araucanía"
$ws.Cells.Item(17, 4).Value = 44494
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 100112026
$ws.Cells.Item(17, 7).Value = "Haba"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 100
$ws.Cells.Item(17, 11).Value = 9000
$ws.Cells.Item(17, 12).Value = 9000
$ws.Cells.Item(17, 13).Value = 9000
$ws.Cells.Item(17, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(17, 15).Value = "Región Metropolitana"
$ws.Cells.Item(17, 16).Value = 360
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = "Hortaliza"
